$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B4"   = 7.724000000000001
    "B6"   = 6.644
    "B7"   = 5.281
    "C7"   = -13.098
    "B8"   = 6.499000000000001
    "C11"  = -12.767
    "C12"  = -11.131
    "C15"  = -13.381
    "B16"  = 5.425000000000001
    "B20"  = 8.099999999999998
    "C20"  = -11.992
    "B21"  = 9.02
    "C21"  = -12.12
    "C22"  = -12.961
    "C23"  = -12.731
    "B28"  = 6.272
    "B29"  = 5.255
    "C29"  = -11.358
    "B30"  = 5.665
    "B32"  = 6.973999999999999
    "C34"  = -12.715
    "B40"  = 9.204000000000001
    "C42"  = -11.999
    "C43"  = -13.688
    "C44"  = -13.599
    "C45"  = -13.339
    "B46"  = 6.063
    "C46"  = -13.953
    "C50"  = -13.938
    "B51"  = 5.059
    "C51"  = -12.093
    "B52"  = 5.880000000000001
    "B57"  = 5.915
    "C57"  = -14.252
    "B59"  = 5.294
    "B62"  = 5.902
    "C65"  = -12.45
    "B66"  = 4.961
    "C66"  = -10.897
    "C67"  = -11.3
    "B73"  = 7.444
    "B74"  = 8.944000000000001
    "B77"  = 6.233000000000001
    "C79"  = -12.204
    "C84"  = -13.643
    "C87"  = -13.691
    "B92"  = 5.761
    "C92"  = -10.898
    "C97"  = -12.751
    "B100" = 6.686
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
